# feat: add 2022-Q4 data
# -----------------------------------------------------------------------
# The workbook tracks quarterly fund-holding snapshots, one worksheet per
# quarter, plus a "总计" (totals) summary sheet. This change adds a new
# "2022-Q4" worksheet (built from a copy of the "2022-Q3" worksheet, since
# both quarters share the same single fund holding/layout) with the Q4
# figures, and records the new quarter in the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q4" sheet by duplicating "2022-Q3" (same layout /
#    formatting), inserting the copy immediately before it so the tab
#    order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# Update the Q4 figures (fund size / stock position / position share /
# held market value) on the new sheet - these are stored as text, same
# as the rest of the sheet, so force a text format before assigning.
$wsQ4.Range("D2:G2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "4.35"
$wsQ4.Range("E2").Value = "98.99"
$wsQ4.Range("F2").Value = "5.30"
$wsQ4.Range("G2").Value = "0.2306"
# Position rank is a real number.
$wsQ4.Range("H2").Value = 3

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert the new 2022-Q4 row at the
#    top of the data and shift the existing rows down. Copy the last
#    data row's formatting into the newly-used row 5 first so the new
#    row matches the sheet's existing look, then fill in all four rows.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("A4:D4").Copy()
$wsTotal.Range("A5:D5").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.23

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.15

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q2"
$wsTotal.Range("C4").Value = 4
$wsTotal.Range("D4").Value = 0.34

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2022-Q1"
$wsTotal.Range("C5").Value = 1
$wsTotal.Range("D5").Value = 0.15
